# Update odds data cells per the 2024-11-16 FlashScore refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65

# Row 3
$ws.Range("G3").Value = 1.53
$ws.Range("H3").Value = 3.75
$ws.Range("J3").Value = 2.2
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("AC3").Value = 7.5
$ws.Range("AD3").Value = 7.5
$ws.Range("AM3").Value = 67
$ws.Range("AT3").Value = 2.5

# Row 4
$ws.Range("G4").Value = 3.6
$ws.Range("I4").Value = 2.05
$ws.Range("J4").Value = 4.33
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 2.75
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 8.5
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.7
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.83
$ws.Range("X4").Value = 17
$ws.Range("AC4").Value = 8.5
$ws.Range("AG4").Value = 301
$ws.Range("AH4").Value = 7
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 67
$ws.Range("BA4").Value = 67

# Row 5
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 3.8
$ws.Range("J5").Value = 2.63
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.75
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.93
$ws.Range("AA5").Value = 17
$ws.Range("AE5").Value = 13
$ws.Range("AH5").Value = 11
$ws.Range("AI5").Value = 19
$ws.Range("AM5").Value = 34
$ws.Range("AQ5").Value = 41

# Row 9
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 3.7
$ws.Range("J9").Value = 2.75
$ws.Range("K9").Value = 2.05
$ws.Range("L9").Value = 4.5
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 2.75
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 1.62
$ws.Range("S9").Value = 1.5
$ws.Range("T9").Value = 2.5
$ws.Range("U9").Value = 1.91
$ws.Range("V9").Value = 1.8
$ws.Range("W9").Value = 6.5
$ws.Range("X9").Value = 9
$ws.Range("Z9").Value = 17
$ws.Range("AD9").Value = 6.5
$ws.Range("AE9").Value = 17
$ws.Range("AF9").Value = 51
$ws.Range("AH9").Value = 9.5
$ws.Range("AI9").Value = 19
$ws.Range("AJ9").Value = 13
$ws.Range("AL9").Value = 34
$ws.Range("AN9").Value = 4
$ws.Range("AO9").Value = 11
$ws.Range("AQ9").Value = 41
$ws.Range("AR9").Value = 67
$ws.Range("AT9").Value = 2.5
$ws.Range("AW9").Value = 5.5
$ws.Range("AX9").Value = 21
$ws.Range("BA9").Value = 101
$ws.Range("BB9").Value = 251

# Row 13
$ws.Range("G13").Value = 2.45
$ws.Range("H13").Value = 3.2
$ws.Range("I13").Value = 2.7
$ws.Range("J13").Value = 3.05
$ws.Range("K13").Value = 2.1
$ws.Range("L13").Value = 3.2
$ws.Range("N13").Value = 6.85
$ws.Range("P13").Value = 2.8
$ws.Range("S13").Value = 1.39
$ws.Range("T13").Value = 2.57
$ws.Range("W13").Value = 7.6
$ws.Range("X13").Value = 11.75
$ws.Range("Y13").Value = 9.5
$ws.Range("Z13").Value = 26
$ws.Range("AA13").Value = 21
$ws.Range("AB13").Value = 32
$ws.Range("AC13").Value = 9
$ws.Range("AD13").Value = 6.2
$ws.Range("AH13").Value = 8.25
$ws.Range("AI13").Value = 13.5
$ws.Range("AJ13").Value = 10
$ws.Range("AK13").Value = 30
$ws.Range("AL13").Value = 23
$ws.Range("AM13").Value = 35
$ws.Range("AN13").Value = 4.35
$ws.Range("AO13").Value = 13
$ws.Range("AP13").Value = 21
$ws.Range("AQ13").Value = 55
$ws.Range("AR13").Value = 90
$ws.Range("AT13").Value = 2.55
$ws.Range("AU13").Value = 7
$ws.Range("AW13").Value = 4.55
$ws.Range("AX13").Value = 14
$ws.Range("AY13").Value = 21
$ws.Range("AZ13").Value = 60
$ws.Range("BA13").Value = 90
$ws.Range("BB13").Value = 250
